$wb = $excel.ActiveWorkbook

$wsSrc = $wb.Worksheets.Item("Week 8")
$wsDst = $wb.Worksheets.Item("Week 9")

# Bring over the date/time number formatting (styles) used by the existing
# entry rows on "Week 8" so the three new rows render as Date / Time / Time.
$wsSrc.Range("A2:C3").Copy()
$wsDst.Range("A2:C4").PasteSpecial(-4122)

# Match the existing row height used by the other data rows.
$wsDst.Rows.Item(2).RowHeight = 18
$wsDst.Rows.Item(3).RowHeight = 18
$wsDst.Rows.Item(4).RowHeight = 18

# New timesheet entries for "Week 9".
$wsDst.Cells.Item(2,1).Value = 43528
$wsDst.Cells.Item(2,2).Value = 0.66666666666666663
$wsDst.Cells.Item(2,3).Value = 0.75
$wsDst.Cells.Item(2,4).Value = "Worked on CSS "
$wsDst.Cells.Item(2,5).Value = 2

$wsDst.Cells.Item(3,1).Value = 43532
$wsDst.Cells.Item(3,2).Value = 0.77083333333333337
$wsDst.Cells.Item(3,3).Value = 0.83333333333333337
$wsDst.Cells.Item(3,4).Value = "Worked on CSS "
$wsDst.Cells.Item(3,5).Value = 1.5

$wsDst.Cells.Item(4,1).Value = 43533
$wsDst.Cells.Item(4,2).Value = 0.5
$wsDst.Cells.Item(4,3).Value = 0.5625
$wsDst.Cells.Item(4,4).Value = "Worked on CSS "
$wsDst.Cells.Item(4,5).Value = 1.5

# Move the active tab / selection from "Week 8" to "Week 9", and place the
# cursor on A5 as in the final saved file.
$wsDst.Select()
$wsDst.Range("A5").Select()
$excel.ActiveWindow.ScrollRow = 11
$excel.ActiveWindow.ScrollColumn = 1
